$d = $word.ActiveDocument
$p1 = $d.Paragraphs.Item(6)
$p2 = $d.Paragraphs.Item(7)
$r = $d.Range($p1.Range.Start, $p2.Range.End)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p><w:pPr><w:pStyle w:val="Normal"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:b w:val="1"/><w:bCs w:val="1"/><w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>5)</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:b w:val="1"/><w:bCs w:val="1"/><w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Download </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>mysql-connector-java-5.0.8-bin.jar and add it to Tomcat Lib directory</w:t></w:r></w:p><w:p w14:noSpellErr="1"><w:pPr><w:pStyle w:val="Normal"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:b w:val="1"/><w:bCs w:val="1"/><w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>6</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:b w:val="1"/><w:bCs w:val="1"/><w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">) </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Setup </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>SSL on Tomcat by following the steps in the how-to guide that exists under wiki &quot;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>How-do-I-set-up-SSL-on-Tomcat.pdf</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>&quot;</w:t></w:r></w:p><w:p w14:noSpellErr="1"><w:pPr><w:pStyle w:val="Normal"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:b w:val="1"/><w:bCs w:val="1"/><w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>7</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:b w:val="1"/><w:bCs w:val="1"/><w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>)</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:b w:val="1"/><w:bCs w:val="1"/><w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Next, we need to update the</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> Tomcat </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>server.xml file to connect to the MySQL</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> realm database</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> that we just created. we need to</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> comment out the UserDatabase realm section. we will create a JDBC</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> realm entry containing all of the correct information, as shown here:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)
Write-Host "done"
